# Updated cryptos list values (price + volume(1h)) per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'69.211.43"
$ws.Range("E2").Value = "  -0.45%  "

$ws.Range("D3").Value = "'3.814.91"
$ws.Range("E3").Value = "  +1.17%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").Value = "'602.27"
$ws.Range("E5").Value = "  -0.67%  "

$ws.Range("D6").Value = "'164.36"
$ws.Range("E6").Value = "  -3.44%  "

$ws.Range("D7").Value = "'3.815.18"
$ws.Range("E7").Value = "  +1.20%  "

$ws.Range("E8").Value = "  -0.02%  "

$ws.Range("E9").Value = "  -0.76%  "

$ws.Range("D10").Value = "'0.171"
$ws.Range("E10").Value = "  +1.02%  "

$ws.Range("D11").Value = "'6.32"
$ws.Range("E11").Value = "  -1.37%  "

$ws.Range("E12").Value = "  -0.37%  "

$ws.Range("D13").Value = "'37.33"
$ws.Range("E13").Value = "  -3.33%  "

$ws.Range("D14").Value = "'0.0000247"

$ws.Range("D15").Value = "'4.452.43"
$ws.Range("E15").Value = "  +1.26%  "

$ws.Range("D16").Value = "'3.805.02"
$ws.Range("E16").Value = "  +1.11%  "

$ws.Range("D17").Value = "'69.318.30"
$ws.Range("E17").Value = "  -0.16%  "

$ws.Range("D18").Value = "'7.45"
$ws.Range("E18").Value = "  +1.53%  "

$ws.Range("E19").Value = "  -0.21%  "

$ws.Range("D20").Value = "'17.36"
$ws.Range("E20").Value = "  +0.92%  "

$ws.Range("D21").Value = "'11.38"
$ws.Range("E21").Value = "  +4.25%  "

$ws.Range("D22").Value = "'489.85"
$ws.Range("E22").Value = "  -1.55%  "

$ws.Range("D23").Value = "'0.724"
$ws.Range("E23").Value = "  -1.12%  "

$ws.Range("D24").Value = "'0.0000157"
$ws.Range("E24").Value = "  +1.12%  "

$ws.Range("D25").Value = "'84.89"
$ws.Range("E25").Value = "  -0.77%  "

$ws.Range("D26").Value = "'2.28"
$ws.Range("E26").Value = "  -3.13%  "

$ws.Range("D27").Value = "'12.28"
$ws.Range("E27").Value = "  -1.21%  "

$ws.Range("D28").Value = "'10.08"
$ws.Range("E28").Value = "  -2.52%  "

$ws.Range("E29").Value = "  +0.05%  "

$ws.Range("E30").Value = "  -0.96%  "

$ws.Range("E31").Value = "  +0.68%  "

$ws.Range("D32").Value = "'2.41"
$ws.Range("E32").Value = "  -4.76%  "

$ws.Range("D33").Value = "'3.956.55"
$ws.Range("E33").Value = "  +1.08%  "

$ws.Range("D34").Value = "'32.06"
$ws.Range("E34").Value = "  -0.12%  "

$ws.Range("D35").Value = "'3.757.83"
$ws.Range("E35").Value = "  +1.55%  "

$ws.Range("E36").Value = "  -1.64%  "

$ws.Range("D37").Value = "'0.140"
$ws.Range("E37").Value = "  +4.88%  "

$ws.Range("E38").Value = "  +0.40%  "

$ws.Range("D39").Value = "'5.94"
$ws.Range("E39").Value = "  +0.04%  "

$ws.Range("D40").Value = "'0.999"
$ws.Range("E40").Value = "  +0.05%  "

$ws.Range("E41").Value = "  -1.51%  "

$ws.Range("D42").Value = "'3.05"
$ws.Range("E42").Value = "  -0.10%  "

$ws.Range("E43").Value = "  -0.07%  "

$ws.Range("D44").Value = "'48.67"
$ws.Range("E44").Value = "  -0.12%  "

$ws.Range("D45").Value = "'424.91"
$ws.Range("E45").Value = "  -3.69%  "

$ws.Range("E46").Value = "  -0.01%  "

$ws.Range("D47").Value = "'8.40"
$ws.Range("E47").Value = "  -1.35%  "

$ws.Range("D48").Value = "'2.834.08"
$ws.Range("E48").Value = "  +0.52%  "

$ws.Range("D49").Value = "'141.41"
$ws.Range("E49").Value = "  +0.31%  "

$ws.Range("D50").Value = "'39.70"
$ws.Range("E50").Value = "  -2.52%  "

$ws.Range("D51").Value = "'0.0351"
$ws.Range("E51").Value = "  -1.48%  "
